$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.459.09"
$ws.Range("E2").Value = "  +0.76%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.612.52"
$ws.Range("E3").Value = "  +1.28%  "

# Row 4
$ws.Range("E4").Value = "  -0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.76"
$ws.Range("E5").Value = "  -0.67%  "

# Row 6
$ws.Range("E6").Value = "  -0.72%  "

# Row 8
$ws.Range("E8").Value = "  -0.50%  "

# Row 9
$ws.Range("E9").Value = "  -0.01%  "

# Row 10
$ws.Range("E10").Value = "  +1.46%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0848"

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.837.81"
$ws.Range("E12").Value = "  +1.18%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.610.17"
$ws.Range("E13").Value = "  +1.10%  "

# Row 14
$ws.Range("E14").Value = "  -0.06%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.509"
$ws.Range("E15").Value = "  +0.00%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.58"
$ws.Range("E16").Value = "  -0.33%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "235.28"
$ws.Range("E17").Value = "  +9.42%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.444.98"
$ws.Range("E18").Value = "  +0.76%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0725"
$ws.Range("E19").Value = "  +0.14%  "

# Row 20
$ws.Range("E20").Value = "  +4.01%  "

# Row 21
$ws.Range("E21").Value = "  +0.06%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.28"
$ws.Range("E22").Value = "  -0.14%  "

# Row 23
$ws.Range("E23").Value = "  +4.28%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.03"
$ws.Range("E24").Value = "  -0.19%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.73"
$ws.Range("E25").Value = "  +1.40%  "

# Row 26
$ws.Range("E26").Value = "  -0.09%  "

# Row 27
$ws.Range("E27").Value = "  +0.29%  "

# Row 28
$ws.Range("E28").Value = "  +0.24%  "

# Row 29
$ws.Range("E29").Value = "  +2.27%  "

# Row 30
$ws.Range("E30").Value = "  +1.05%  "

# Row 31
$ws.Range("E31").Value = "  -0.47%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.496.36"
$ws.Range("E32").Value = "  +6.17%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.24"
$ws.Range("E33").Value = "  +1.32%  "

# Row 34
$ws.Range("E34").Value = "  -0.98%  "

# Row 35
$ws.Range("E35").Value = "  -0.45%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.49"
$ws.Range("E36").Value = "  +1.82%  "

# Row 37
$ws.Range("E37").Value = "  -2.72%  "

# Row 38
$ws.Range("E38").Value = "  -0.09%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.827"
$ws.Range("E39").Value = "  +0.61%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.78"
$ws.Range("E40").Value = "  +0.23%  "

# Row 41
$ws.Range("E41").Value = "  -0.05%  "

# Row 42
$ws.Range("E42").Value = "  +1.15%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.929"
$ws.Range("E43").Value = "  -3.08%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.750.47"
$ws.Range("E44").Value = "  +1.32%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.760"
$ws.Range("E45").Value = "  -0.02%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "61.34"
$ws.Range("E46").Value = "  +0.80%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "89.77"
$ws.Range("E47").Value = "  +3.10%  "

# Row 48
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.49"
$ws.Range("E48").Value = "  -0.19%  "

# Row 49
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0501"
$ws.Range("E49").Value = "  +0.03%  "

# Row 50
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0961"
$ws.Range("E50").Value = "  +0.96%  "

# Row 51
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.47"
$ws.Range("E51").Value = "  +1.32%  "
